$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.242.05"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.992.03"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB (value looks numeric, force text)
$ws.Range("D5").Value = "'560.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6 - Solana (value looks numeric, force text)
$ws.Range("D6").Value = "'137.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.46%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - XRP (value looks numeric, force text)
$ws.Range("D8").Value = "'0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.981.67"
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - Dogecoin (value looks numeric, force text)
$ws.Range("D10").Value = "'0.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.37%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +4.83%  "

# Row 12 - Cardano (value looks numeric, force text)
$ws.Range("D12").Value = "'0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.87%  "

# Row 13 - ShibaInu (value looks numeric, force text)
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.89%  "

# Row 14 - Avalanche (value looks numeric, force text)
$ws.Range("D14").Value = "'33.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.91%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.40%  "

# Row 16 - now WrappedliquidstakedEther2.0 (was Polkadot)
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.485.91"
$ws.Range("E16").Value = "  +0.41%  "

# Row 17 - now Polkadot (was WrappedliquidstakedEther2.0) (value looks numeric, force text)
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'7.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.30%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.989.57"
$ws.Range("E18").Value = "  +0.45%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "59.242.47"
$ws.Range("E19").Value = "  +1.68%  "

# Row 20 - BitcoinCash (value looks numeric, force text)
$ws.Range("D20").Value = "'428.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "

# Row 21 - Chainlink (value looks numeric, force text)
$ws.Range("D21").Value = "'13.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.87%  "

# Row 22 - Polygon (value looks numeric, force text)
$ws.Range("D22").Value = "'0.720"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.95%  "

# Row 23 - InternetComputer(DFINITY) (value looks numeric, force text)
$ws.Range("D23").Value = "'13.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.57%  "

# Row 24 - Uniswap
$ws.Range("E24").Value = "  +1.01%  "

# Row 25 - Litecoin (value looks numeric, force text)
$ws.Range("D25").Value = "'80.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.23%  "

# Row 27 - ImmutableX (value looks numeric, force text)
$ws.Range("D27").Value = "'2.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.65%  "

# Row 28 - FirstDigitalUSD
$ws.Range("E28").Value = "  +0.13%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.75%  "

# Row 30 - RenderToken (value looks numeric, force text)
$ws.Range("D30").Value = "'7.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "

# Row 31 - EthereumClassic (value looks numeric, force text)
$ws.Range("D31").Value = "'25.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.44%  "

# Row 32 - NEARProtocol (value looks numeric, force text)
$ws.Range("D32").Value = "'6.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.13%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +1.17%  "

# Row 34 - Mantle (value looks numeric, force text)
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.89%  "

# Row 35 - Filecoin (value looks numeric, force text)
$ws.Range("D35").Value = "'5.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.81%  "

# Row 36 - PEPE
$ws.Range("D36").Value = "0.0₃0757"
$ws.Range("E36").Value = "  +7.60%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -2.33%  "

# Row 38 - OKB (value looks numeric, force text)
$ws.Range("D38").Value = "'48.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "

# Row 39 - Cosmos (value looks numeric, force text)
$ws.Range("D39").Value = "'8.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "

# Row 40 - dogwifhat
$ws.Range("E40").Value = "  +4.46%  "

# Row 41 - Bittensor (value looks numeric, force text)
$ws.Range("D41").Value = "'405.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.30%  "

# Row 42 - now Maker (was VeChain)
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.772.39"
$ws.Range("E42").Value = "  +1.89%  "

# Row 43 - now VeChain (was Maker) (value looks numeric, force text)
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0352"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.40%  "

# Row 44 - Kaspa
$ws.Range("E44").Value = "  -1.88%  "

# Row 45 - TheGraph (value looks numeric, force text)
$ws.Range("D45").Value = "'0.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.00%  "

# Row 47 - Monero (value looks numeric, force text)
$ws.Range("D47").Value = "'123.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "

# Row 48 - Arweave (value looks numeric, force text)
$ws.Range("D48").Value = "'34.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +19.17%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -0.54%  "

# Row 50 - Fetch.AI (value looks numeric, force text)
$ws.Range("D50").Value = "'2.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.76%  "

# Row 51 - InjectiveProtocol (value looks numeric, force text)
$ws.Range("D51").Value = "'23.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.06%  "
